$d = $word.ActiveDocument

$text1 = "1. the validation code will trigger whenever there is a Insert operation sent to the sql database. In this specific case it will check that the object being created has a valid extension by matching the extension with a regular expression pattern."
$text2 = '2. First a function is defined that returns a test product object. Then the test "image url" is run. This test makes 2 lists of url extensions, one with correct extensions, and the other with incorrect. It then asserts that a product created with each one of the good names is valid, and then asserts that a product created using each of the bad names is invalid.'
$text3 = "3. A test fixture specifies the initial contents of a model under a test. The fixtures directive loads the corresponding fixture data into the corresponding database table before each test is run."
$text4 = "4. Overdue_Book"
$text5 = "5. The three databases are development, test, and production"

$newText = $text1 + "`r" + $text2 + "`r" + $text3 + "`r" + $text4 + "`r" + $text5

$count = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($count)
$r = $p.Range
$r.InsertBefore($newText)
